$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11
$ws.Range("F11").Value = 415
$ws.Range("H11").Value = 480

# Row 12
$ws.Range("F12").Value = 672
$ws.Range("H12").Value = 758

# Row 23
$ws.Range("F23").Value = 175
$ws.Range("H23").Value = 227

# Row 24
$ws.Range("F24").Value = 247
$ws.Range("H24").Value = 277

# Row 27
$ws.Range("F27").Value = 315
$ws.Range("H27").Value = 397

# Row 31
$ws.Range("F31").Value = 54
$ws.Range("H31").Value = 81

# Row 33
$ws.Range("E33").Value = 376
$ws.Range("F33").Value = 244
$ws.Range("H33").Value = 335

# Row 35
$ws.Range("E35").Value = 198
$ws.Range("F35").Value = 168
$ws.Range("H35").Value = 195

# Row 39
$ws.Range("F39").Value = 149
$ws.Range("H39").Value = 200

# Row 44
$ws.Range("F44").Value = 316
$ws.Range("H44").Value = 384

# Row 46
$ws.Range("F46").Value = 327
$ws.Range("H46").Value = 391
